$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update "想去人数" (column F) values
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F17").Value = 1537
$wsExhibit.Range("F22").Value = 8289
$wsExhibit.Range("F27").Value = 1282
$wsExhibit.Range("F31").Value = 6620
$wsExhibit.Range("F37").Value = 5670

# Sheet "全部类型" (sheet4): same events duplicated, update matching rows
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F17").Value = 1537
$wsAll.Range("F23").Value = 8289
$wsAll.Range("F28").Value = 1282
$wsAll.Range("F34").Value = 6620
$wsAll.Range("F40").Value = 5670
